$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.627.81"
$ws.Range("E2").Value = "  -1.92%  "

$ws.Range("D3").Value = "2.348.65"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("D5").Value = "'323.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.80%  "

$ws.Range("D6").Value = "'101.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.42%  "

$ws.Range("D7").Value = "'0.636"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.08%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.17%  "

$ws.Range("D10").Value = "'39.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.40%  "

$ws.Range("D11").Value = "'0.0918"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.29%  "

$ws.Range("D12").Value = "'8.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.31%  "

$ws.Range("D13").Value = "'0.995"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.55%  "

$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("D15").Value = "'16.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.16%  "

$ws.Range("D16").Value = "2.709.02"
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("D17").Value = "2.359.56"
$ws.Range("E17").Value = "  -2.91%  "

$ws.Range("D18").Value = "'7.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.59%  "

$ws.Range("D19").Value = "42.617.54"
$ws.Range("E19").Value = "  -1.92%  "

$ws.Range("E20").Value = "  -2.76%  "

$ws.Range("D21").Value = "'75.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("D22").Value = "'3.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.05%  "

$ws.Range("D23").Value = "'264.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.51%  "

$ws.Range("E24").Value = "  -9.53%  "

$ws.Range("E25").Value = "  +8.26%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "'11.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.90%  "

$ws.Range("D28").Value = "'22.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.76%  "

$ws.Range("E29").Value = "  -2.18%  "

$ws.Range("D30").Value = "'175.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.20%  "

$ws.Range("E31").Value = "  -3.03%  "

$ws.Range("D32").Value = "'0.0898"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.67%  "

$ws.Range("D33").Value = "'35.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.63%  "

$ws.Range("E34").Value = "  +0.43%  "

$ws.Range("E35").Value = "  -0.92%  "

$ws.Range("D36").Value = "'4.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.94%  "

$ws.Range("D37").Value = "'0.0356"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.68%  "

$ws.Range("E38").Value = "  +2.27%  "

$ws.Range("D39").Value = "'2.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.73%  "

$ws.Range("D40").Value = "'3.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.54%  "

$ws.Range("D41").Value = "'1.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("D42").Value = "'0.237"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.42%  "

$ws.Range("D43").Value = "'69.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.25%  "

$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").Value = "'119.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.40%  "

$ws.Range("D46").Value = "'91.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +20.10%  "

$ws.Range("D47").Value = "'11.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.18%  "

$ws.Range("E48").Value = "  -2.15%  "

$ws.Range("D49").Value = "'9.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.50%  "

$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("E51").Value = "  -4.51%  "
